# Auto-generated edit script applying the Sophia_Profits.xlsx market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()
$ws.Cells.Item(12, 8).Value = 149.42857
$ws.Cells.Item(12, 9).Value = 181.8
$ws.Cells.Item(12, 10).Value = 68.5
$ws.Cells.Item(12, 11).Value = 181.8
$ws.Cells.Item(12, 12).Value = 68.5
$ws.Cells.Item(12, 13).Value = -11.80000000000001
$ws.Cells.Item(12, 14).Value = -408.5
$ws.Cells.Item(28, 8).Value = 2081.8
$ws.Cells.Item(28, 9).Value = 3421.3333
$ws.Cells.Item(28, 10).Value = 72.5
$ws.Cells.Item(28, 11).Value = 3421.3333
$ws.Cells.Item(28, 12).Value = 72.5
$ws.Cells.Item(28, 13).Value = -2936.3333
$ws.Cells.Item(28, 14).Value = -1042.5
$ws.Cells.Item(80, 8).Value = 7078
$ws.Cells.Item(80, 9).Value = 6733.3335
$ws.Cells.Item(80, 10).Value = 7595
$ws.Cells.Item(80, 11).Value = 20200.0005
$ws.Cells.Item(80, 12).Value = 22785
$ws.Cells.Item(80, 13).Value = -19202.0005
$ws.Cells.Item(80, 14).Value = -24781
$ws.Cells.Item(83, 8).Value = 7078
$ws.Cells.Item(83, 9).Value = 6733.3335
$ws.Cells.Item(83, 10).Value = 7595
$ws.Cells.Item(83, 11).Value = 60600.0015
$ws.Cells.Item(83, 12).Value = 68355
$ws.Cells.Item(83, 13).Value = -55608.0015
$ws.Cells.Item(83, 14).Value = -78339
$ws.Cells.Item(92, 8).Value = 822.4286
$ws.Cells.Item(92, 9).Value = 747.2308
$ws.Cells.Item(92, 11).Value = 747.2308
$ws.Cells.Item(92, 13).Value = 500.7692
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 2661.625
$ws.Cells.Item(137, 10).Value = 1998.5
$ws.Cells.Item(137, 12).Value = 5995.5
$ws.Cells.Item(137, 14).Value = -11095.5
$ws.Cells.Item(138, 8).Value = 4048.5557
$ws.Cells.Item(138, 10).Value = 5011.5264
$ws.Cells.Item(138, 12).Value = 15034.5792
$ws.Cells.Item(138, 14).Value = -25314.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5135.3228
$ws.Cells.Item(32, 9).Value = 5135.3228
$ws.Cells.Item(32, 11).Value = 5135.3228
$ws.Cells.Item(32, 13).Value = -4848.3228
$ws.Cells.Item(74, 8).Value = 10148.75
$ws.Cells.Item(74, 9).Value = 11415.08
$ws.Cells.Item(74, 10).Value = 5626.143
$ws.Cells.Item(74, 11).Value = 11415.08
$ws.Cells.Item(74, 12).Value = 5626.143
$ws.Cells.Item(74, 13).Value = -10541.08
$ws.Cells.Item(74, 14).Value = -7374.143
$ws.Cells.Item(77, 8).Value = 10148.75
$ws.Cells.Item(77, 9).Value = 11415.08
$ws.Cells.Item(77, 10).Value = 5626.143
$ws.Cells.Item(77, 11).Value = 57075.4
$ws.Cells.Item(77, 12).Value = 28130.715
$ws.Cells.Item(77, 13).Value = -52707.4
$ws.Cells.Item(77, 14).Value = -36866.715
$ws.Cells.Item(110, 8).Value = 3360.4167
$ws.Cells.Item(110, 9).Value = 1488.3
$ws.Cells.Item(110, 10).Value = 12721
$ws.Cells.Item(110, 11).Value = 1488.3
$ws.Cells.Item(110, 12).Value = 12721
$ws.Cells.Item(110, 13).Value = 556.7
$ws.Cells.Item(110, 14).Value = -16811
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4108.6924
$ws.Cells.Item(20, 9).Value = 1046.3334
$ws.Cells.Item(20, 11).Value = 1046.3334
$ws.Cells.Item(20, 13).Value = -799.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2225
$ws.Cells.Item(31, 9).Value = 2366.6667
$ws.Cells.Item(31, 11).Value = 2366.6667
$ws.Cells.Item(31, 13).Value = -2071.6667
$ws.Cells.Item(34, 8).Value = 2225
$ws.Cells.Item(34, 9).Value = 2366.6667
$ws.Cells.Item(34, 11).Value = 2366.6667
$ws.Cells.Item(34, 13).Value = -2164.6667
$ws.Cells.Item(41, 8).Value = 22990
$ws.Cells.Item(41, 10).Value = 22990
$ws.Cells.Item(41, 12).Value = 22990
$ws.Cells.Item(41, 14).Value = -23846
$ws.Cells.Item(58, 8).Value = 4673.778
$ws.Cells.Item(58, 9).Value = 4673.778
$ws.Cells.Item(58, 11).Value = 4673.778
$ws.Cells.Item(58, 13).Value = -4470.778
$ws.Cells.Item(136, 8).Value = 4673.778
$ws.Cells.Item(136, 9).Value = 4673.778
$ws.Cells.Item(136, 11).Value = 14021.334
$ws.Cells.Item(136, 13).Value = -11471.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 160.66667
$ws.Cells.Item(2, 9).Value = 92.25
$ws.Cells.Item(2, 10).Value = 297.5
$ws.Cells.Item(2, 11).Value = 553.5
$ws.Cells.Item(2, 12).Value = 1785
$ws.Cells.Item(2, 13).Value = -440.5
$ws.Cells.Item(2, 14).Value = -2011
$ws.Cells.Item(12, 8).Value = 239.47058
$ws.Cells.Item(12, 9).Value = 357.42856
$ws.Cells.Item(12, 10).Value = 156.9
$ws.Cells.Item(12, 11).Value = 1072.28568
$ws.Cells.Item(12, 12).Value = 470.7
$ws.Cells.Item(12, 13).Value = -899.28568
$ws.Cells.Item(12, 14).Value = -816.7
$ws.Cells.Item(14, 8).Value = 520.2727
$ws.Cells.Item(14, 9).Value = 520.2727
$ws.Cells.Item(14, 11).Value = 1560.8181
$ws.Cells.Item(14, 13).Value = -1387.8181
$ws.Cells.Item(55, 8).Value = 11660
$ws.Cells.Item(55, 10).Value = 11660
$ws.Cells.Item(55, 12).Value = 34980
$ws.Cells.Item(55, 14).Value = -35334
$ws.Cells.Item(87, 8).Value = 7244.625
$ws.Cells.Item(87, 9).Value = 6657.6665
$ws.Cells.Item(87, 10).Value = 9005.5
$ws.Cells.Item(87, 11).Value = 19972.9995
$ws.Cells.Item(87, 12).Value = 27016.5
$ws.Cells.Item(87, 13).Value = -18724.9995
$ws.Cells.Item(87, 14).Value = -29512.5
$ws.Cells.Item(90, 8).Value = 7244.625
$ws.Cells.Item(90, 9).Value = 6657.6665
$ws.Cells.Item(90, 10).Value = 9005.5
$ws.Cells.Item(90, 11).Value = 59918.9985
$ws.Cells.Item(90, 12).Value = 81049.5
$ws.Cells.Item(90, 13).Value = -53678.9985
$ws.Cells.Item(90, 14).Value = -93529.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3872.625
$ws.Cells.Item(102, 9).Value = 3997.2856
$ws.Cells.Item(102, 11).Value = 3997.2856
$ws.Cells.Item(102, 13).Value = -2375.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(76, 8).Value = 69884
$ws.Cells.Item(76, 10).Value = 69884
$ws.Cells.Item(76, 12).Value = 69884
$ws.Cells.Item(76, 14).Value = -70560
$ws.Cells.Item(79, 8).Value = 69884
$ws.Cells.Item(79, 10).Value = 69884
$ws.Cells.Item(79, 12).Value = 69884
$ws.Cells.Item(79, 14).Value = -72224
$ws.Cells.Item(106, 8).Value = 48599.25
$ws.Cells.Item(106, 10).Value = 48599.25
$ws.Cells.Item(106, 12).Value = 48599.25
$ws.Cells.Item(106, 14).Value = -51123.25
$ws.Cells.Item(136, 8).Value = 3498.5
$ws.Cells.Item(136, 9).Value = 3498.5
$ws.Cells.Item(136, 11).Value = 10495.5
$ws.Cells.Item(136, 13).Value = -7945.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5340
$ws.Cells.Item(81, 9).Value = 2063.3333
$ws.Cells.Item(81, 11).Value = 4126.6666
$ws.Cells.Item(81, 13).Value = -3065.6666
$ws.Cells.Item(84, 8).Value = 5340
$ws.Cells.Item(84, 9).Value = 2063.3333
$ws.Cells.Item(84, 11).Value = 20633.333
$ws.Cells.Item(84, 13).Value = -15329.333
$ws.Cells.Item(123, 8).Value = 100000
$ws.Cells.Item(123, 10).Value = 100000
$ws.Cells.Item(123, 12).Value = 100000
$ws.Cells.Item(123, 14).Value = -109800
$ws.Cells.Item(126, 8).Value = 1524.125
$ws.Cells.Item(126, 10).Value = 1197
$ws.Cells.Item(126, 12).Value = 3591
$ws.Cells.Item(126, 14).Value = -8531
$ws.Cells.Item(130, 8).Value = 93999
$ws.Cells.Item(130, 10).Value = 93999
$ws.Cells.Item(130, 12).Value = 93999
$ws.Cells.Item(130, 14).Value = -104039
$ws.Cells.Item(136, 8).Value = 4718.8
$ws.Cells.Item(136, 9).Value = 4718.8
$ws.Cells.Item(136, 11).Value = 14156.4
$ws.Cells.Item(136, 13).Value = -11606.4
